# Update "Datos actualizados..." timestamp and refresh the provincias_spain
# case-count table (re-sorted by "Casos totales" descending, new values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 13:16"

$data = @(
    @("Madrid", 14597, 3031, 9741, 1825),
    @("Cataluña", 9937, 1274, 8147, 516),
    @("Castilla-La Mancha", 2780, 71, 2446, 263),
    @("Bizkaia/Vizcaya", 1501, 466, 1457, 44),
    @("Valencia/Valencia", 1497, 23, 1413, 61),
    @("Araba/Alava", 1207, 466, 1134, 73),
    @("Navarra", 1197, 23, 1141, 33),
    @("La Rioja", 928, 43, 848, 37),
    @("Aragon", 907, 22, 845, 40),
    @("Alacant/Alicante", 857, 12, 774, 71),
    @("Malaga", 819, 72, 793, 133),
    @("Asturias", 779, 35, 719, 25),
    @("A Coruña", 734, 25, 715, 19),
    @("Pontevedra", 602, 25, 598, 4),
    @("Granada", 579, 72, 558, 133),
    @("Gipuzkoa/Guipuzcoa", 563, 466, 547, 16),
    @("Sevilla", 535, 72, 521, 133),
    @("Cantabria", 510, 12, 484, 14),
    @("Ciudad Real", 505, 8, 457, 40),
    @("Toledo", 501, 22, 451, 28),
    @("Caceres", 485, 3, 447, 35),
    @("Salamanca", 483, 42, 403, 38),
    @("Murcia", 477, 4, 467, 6),
    @("Tenerife", 438, 8, 394, 21),
    @("Albacete", 430, 8, 390, 32),
    @("Valladolid", 410, 24, 369, 17),
    @("Burgos", 392, 41, 327, 24),
    @("Leon", 362, 21, 317, 24),
    @("Zaragoza", 329, 0, 315, 14),
    @("Jaen", 316, 72, 304, 133),
    @("Segovia", 271, 32, 212, 27),
    @("Castello/Castellon", 269, 1, 257, 11),
    @("Guadalajara", 263, 2, 257, 4),
    @("Badajoz", 257, 5, 248, 4),
    @("Mallorca", 210, 18, 194, 12),
    @("Avila", 201, 23, 163, 15),
    @("Cordoba", 191, 72, 185, 133),
    @("Ourense", 189, 25, 186, 3),
    @("Soria", 179, 14, 152, 13),
    @("Cadiz", 178, 72, 175, 133),
    @("Gran Canaria", 158, 8, 156, 21),
    @("Cuenca", 120, 8, 104, 8),
    @("Almeria", 115, 72, 110, 133),
    @("Lugo", 103, 25, 99, 4),
    @("Zamora", 100, 11, 82, 7),
    @("Huelva", 77, 72, 76, 133),
    @("Palencia", 72, 13, 59, 0)
)

$row = 4
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
